# Journal de veille technologique - add "Menu Pause" entry and a link to the
# existing "Matériel" entry; update consultation dates for the Terrain and
# RayCast entries (8 février 2021), matching commit:
# "Ajout d'un menu principal, ainsi qu'un menu de pause"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal de veille")

# Row 14 ("Matériel" source): append the reference URL under the title.
$ws.Range("D14").Value = "Matériel`nhttps://docs.unity3d.com/ScriptReference/Material.html"

# Row 15 ("Terrain" source): consultation date moves from 5 to 8 février 2021.
$ws.Range("E15").Value = "8 février 2021"

# Row 16 ("RayCast" source): consultation date moves from 5 to 8 février 2021.
$ws.Range("E16").Value = "8 février 2021"

# Row 17 was an almost-empty placeholder row; fill it in with the new
# "Menu Pause" video source (order matters for shared-string allocation:
# the date is entered before the title/summary, as in the original edit).
$ws.Range("B17").Value = "Youtube"
$ws.Range("C17").Value = "Réseaux sociaux"
$ws.Range("E17").Value = "9 février 2021"
$ws.Range("D17").Value = "Menu Pause`nhttps://www.youtube.com/watch?v=JivuXdrIHK0"
$ws.Range("F17").Value = "Cette vidéo permet de comprendre les bases de la création d'un menu pause, ainsi que les utilisations de celui-ci. C'est-à-dire le retour au menu principale, et même quitter le jeux."
$ws.Rows.Item(17).RowHeight = 45

# Reflect where the author's cursor ended up after the edit.
$ws.Range("M16").Select()
